$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. New header row: insert WIN / TOP4 / TOP5 / RELEGATION before the
#    existing ExpPoints column, which shifts from C to G.
# ------------------------------------------------------------------

# Push the current C column ("ExpPoints" header + values) out to G
# first, then overwrite C1:F1 with the four new headers.
$ws.Range("C1").Copy()
$ws.Range("G1").PasteSpecial(-4104)   # xlPasteAll (value + formatting)

$headers = @("WIN", "TOP4", "TOP5", "RELEGATION")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 3 + $i   # C, D, E, F
    $ws.Cells.Item(1, $col).Value = $headers[$i]
}

# Copy the header formatting (bold font + thin border + centred
# alignment, style index "1" in the original file) onto the new
# header cells so they match the existing Rank/Team/ExpPoints look.
$ws.Range("B1").Copy()
$ws.Range("C1:G1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("G1").Value = "ExpPoints"

# ------------------------------------------------------------------
# 2. Data rows: new ordering (Aston Villa / Newcastle United /
#    Brighton & Hove Albion / AFC Bournemouth re-ranked) together with
#    refreshed ExpPoints values, now living in column G. Columns C-F
#    are left blank placeholders for the upcoming Monte Carlo
#    percentages (WIN / TOP4 / TOP5 / RELEGATION).
# ------------------------------------------------------------------

$teams = @(
    @("Arsenal", 79.27864085096618),
    @("Liverpool", 73.35949446806187),
    @("Manchester City", 70.88558983835186),
    @("Crystal Palace", 60.71348164513891),
    @("Chelsea", 59.77805889766427),
    @("Aston Villa", 55.58089203449666),
    @("Newcastle United", 54.72658192009546),
    @("Brighton & Hove Albion", 54.6943194037909),
    @("AFC Bournemouth", 54.60953215016666),
    @("Tottenham Hotspur", 54.21988475446583),
    @("Manchester United", 51.60137344874742),
    @("Brentford", 50.0591257723337),
    @("Fulham", 44.59697644818363),
    @("Everton", 44.11525500300813),
    @("Sunderland", 42.09535988989686),
    @("Nottingham Forest", 38.95415329593722),
    @("Leeds United", 37.1361237389611),
    @("West Ham United", 36.05463518645836),
    @("Burnley", 33.26020270339857),
    @("Wolverhampton Wanderers", 29.85683734269547)
)

for ($i = 0; $i -lt $teams.Length; $i++) {
    $row = $i + 2
    $team = $teams[$i][0]
    $points = $teams[$i][1]

    $ws.Cells.Item($row, 2).Value = $team          # B: Team
    $ws.Cells.Item($row, 3).Value = ""              # C: WIN (blank)
    $ws.Cells.Item($row, 4).Value = ""              # D: TOP4 (blank)
    $ws.Cells.Item($row, 5).Value = ""              # E: TOP5 (blank)
    $ws.Cells.Item($row, 6).Value = ""              # F: RELEGATION (blank)
    $ws.Cells.Item($row, 7).Value = $points          # G: ExpPoints
}

Write-Output "applied"
